$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "DTaP" -> "DTaP/"  (rows 2-4, column A)
$ws.Range("A2").Value = "DTaP/"
$ws.Range("A3").Value = "DTaP/"
$ws.Range("A4").Value = "DTaP/"

# "DTaP " -> "DTaP-Hib "  (row 5, column A)
$ws.Range("A5").Value = "DTaP-Hib "

# "Hepatitis Bi " -> "Hepatitis B-Hibi "  (row 7, column A)
$ws.Range("A7").Value = "Hepatitis B-Hibi "

# "Hepatitis BPreservative Free PediatricAdolescent" -> "Hepatitis B Preservative Free Pediatric/Adolescent" (rows 17-22, column A)
$ws.Range("A17").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"
$ws.Range("A18").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"
$ws.Range("A19").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"
$ws.Range("A20").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"
$ws.Range("A21").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"
$ws.Range("A22").Value = "Hepatitis B Preservative Free Pediatric/Adolescent"

# "251 dose TipLok 23G" -> "25 x 1 dose TipLok 23G" (row 19, column D)
$ws.Range("D19").Value = "25 x 1 dose TipLok 23G"

# "251 dose TipLok 25G" -> "25 x 1 dose TipLok 25G" (row 20, column D)
$ws.Range("D20").Value = "25 x 1 dose TipLok 25G"

# "11 dose TipLok" -> "1x1 dose TipLok" (row 21, column D)
$ws.Range("D21").Value = "1x1 dose TipLok"

# "Hepatitis B 2 doseAdolescent (11-15)" -> "Hepatitis B 2 dose Adolescent (11-15)" (row 23, column A)
$ws.Range("A23").Value = "Hepatitis B 2 dose Adolescent (11-15)"

# "MMR" -> "MMR/" (row 29, column A)
$ws.Range("A29").Value = "MMR/"

# "Pneumococcal7-valent (Pediatric)" -> "Pneumococcal 7-valent (Pediatric)" (row 30, column A)
$ws.Range("A30").Value = "Pneumococcal 7-valent (Pediatric)"
